$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Corrigiendo permisos del Regente: quitar "formulas" de la lista de permisos
$ws.Range("B3").Value = "CRUD(medicamentos,sucursales)"

# Reflejar el ultimo estado de seleccion de la hoja (B18)
$ws.Range("B18").Select()
